{"js": "// Insert a new bulleted list item right after the paragraph that ends with\n// \"Finished the 2d model for the chassis and 3d model of the propeller\",\n// matching the existing list formatting (ListParagraph style, numId 3).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"Finished the 2d model for the chassis and 3d model of the propeller\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\nconst newText =\n  \"We wanted the design of the chassis to work differently than it does in the final version. We would have liked it to blow air in the direction of the components, but due to limitations of the flexibility of the power cable we opted to have it blow \\u201Cbackwards\\u201D compared to conventional fan design.\";\n\n// insertParagraph after the anchor paragraph inherits its paragraph\n// formatting (pStyle \"ListParagraph\" + numPr numId 3, lang en-GB), matching\n// the diff's new <w:p>.\nconst newPara = anchor.insertParagraph(newText, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Insert a new bulleted list item right after the paragraph that ends with\n# \"Finished the 2d model for the chassis and 3d model of the propeller\",\n# matching the existing list formatting (ListParagraph style, numId 3).\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Finished the 2d model for the chassis and 3d model of the propeller\"\n$leftQuote = [char]0x201C\n$rightQuote = [char]0x201D\n$newText = \"We wanted the design of the chassis to work differently than it does in the final version. We would have liked it to blow air in the direction of the components, but due to limitations of the flexibility of the power cable we opted to have it blow \" + $leftQuote + \"backwards\" + $rightQuote + \" compared to conventional fan design.\"\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq $anchorText) {\n        # InsertParagraphAfter() creates a new paragraph right after this one,\n        # inheriting its paragraph formatting (style + numbering + rPr).\n        $p.Range.InsertParagraphAfter()\n        $newPara = $p.Next()\n        $newPara.Range.Text = $newText\n        break\n    }\n}\n"}
